$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3937
$ws1.Range("F4").Value = 2320
$ws1.Range("F9").Value = 112
$ws1.Range("F10").Value = 22
$ws1.Range("F11").Value = 118
$ws1.Range("F12").Value = 1460
$ws1.Range("F14").Value = 2659

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3937
$ws4.Range("F4").Value = 2320
$ws4.Range("F10").Value = 112
$ws4.Range("F11").Value = 22
$ws4.Range("F12").Value = 118
$ws4.Range("F15").Value = 1460
$ws4.Range("F17").Value = 2659
